$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) values for rows 2-51.
# Values are prefixed with a leading apostrophe so Excel keeps them as
# literal text (matching the original inline-string cells) instead of
# reinterpreting number-like strings (e.g. "1.00", "507.83") as numbers.
# The Style is reset to "Normal" afterwards so no new cell format/style
# is introduced, keeping the cells styled exactly as before.

$ws.Cells.Item(2, 4).Value = "'57.204.73"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  +1.14%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'3.001.40"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -0.29%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  +0.06%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'507.83"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -0.12%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'138.49"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  +1.18%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 4).Value = "'1.00"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  +0.01%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'0.436"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  +0.48%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'7.50"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -1.76%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 4).Value = "'0.110"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  +0.76%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'0.365"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  +3.02%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 4).Value = "'3.534.77"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  +0.28%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'0.130"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  +1.54%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 4).Value = "'26.36"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  +2.43%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'0.0000164"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  +6.29%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'57.293.37"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  +1.19%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'6.20"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  +5.84%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = "'3.021.19"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  +0.40%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 4).Value = "'12.70"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  +1.51%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'7.92"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  +0.61%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'328.56"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  +0.26%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'0.999"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -0.03%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'0.495"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +3.34%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'64.38"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  +2.95%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'0.168"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  +0.86%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'0.998"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -0.40%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 4).Value = "'0.0₃0912"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  -0.95%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 4).Value = "'6.71"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  +2.06%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "'7.41"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  +5.83%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'1.80"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  +1.69%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).Value = "'1.18"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -5.35%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'20.56"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -0.35%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'4.68"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  +3.28%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'154.14"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -1.23%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "'5.87"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  +4.41%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'1.27"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -0.64%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 4).Value = "'24.32"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  +2.52%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'0.0678"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  +0.42%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'3.045.81"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  +0.03%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'37.16"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  +1.55%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  +0.15%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'3.82"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +5.80%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "'2.292.73"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  +1.14%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'0.650"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +0.32%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'1.40"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -0.67%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 4).Value = "'0.977"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -2.57%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'6.00"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  +3.28%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 4).Value = "'0.0238"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  +0.79%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'19.29"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  +0.53%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'1.83"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  -8.20%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 4).Value = "'0.0891"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  +1.79%  "
$ws.Cells.Item(51, 5).Style = "Normal"
